# typo: fix typo in members sheet.
#
# - Sheet3 column C ("description" of the "Third" batch) is filled with the
#   misspelling "Coordinatior" for every member row (C2:C16) -> "Coordinator".
# - Sheet4 row 3 ("Marmik Sharma") has "Technical Secetary" -> "Technical Secretary".
# - The workbook's active sheet moves from Sheet3 to Sheet4, and the
#   remembered selection on each of those sheets changes accordingly.

$wb = $excel.ActiveWorkbook

# --- Fix "Coordinatior" -> "Coordinator" on Sheet3 (all member rows, col C) ---
$ws3 = $wb.Worksheets.Item("Sheet3")
for ($row = 2; $row -le 16; $row++) {
    $cell = $ws3.Cells.Item($row, 3)
    if ($cell.Value2 -eq "Coordinatior") {
        $cell.Value = "Coordinator"
    }
}

# --- Fix "Technical Secetary" -> "Technical Secretary" on Sheet4 ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("C3").Value = "Technical Secretary"

# --- Update remembered selection on Sheet3, then move the active sheet/
#     selection to Sheet4 (matches workbookView activeTab 2 -> 3) ---
$ws3.Range("D20").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("D24").Select() | Out-Null
